# Fix a copy/paste error in the assignment documentation table.
#
# The "Assignment without &" sheet (sheet1 / rId1) had three rows (E4, E13,
# E20, E29) that were copy/pasted from rows describing "(RHS.R), RHS.R"
# style moves, but should actually describe the "(RHS.x), R | MOV R, (LHS.x)"
# (LEA-style) moves. Also the workbook had accidentally been saved with the
# second tab ("Assignment with &") active/selected instead of the first one.

$wb = $excel.ActiveWorkbook

$wsWithout = $wb.Worksheets.Item("Assignment without &")

# Correct the copy/pasted cell text in the "Assignment without &" sheet.
$wsWithout.Range("E4").Value  = "MOV (RHS.A), R | MOV R, (LHS.A)"
$wsWithout.Range("E13").Value = "MOV (RHS.R), R | MOV R, (LHS.A)"
$wsWithout.Range("E20").Value = "MOV (RHS.A), R | MOV R, (LHS.R)"
$wsWithout.Range("E29").Value = "MOV (RHS.A), R | MOV R, (LHS.R)"

# Restore selection: first tab ("Assignment without &") should be the
# active/selected sheet, not the second one.
$wsWithout.Select()
$wsWithout.Activate()
